$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column F: mark the existing bug as Fixed = "Yes"
$ws.Range("F2").Value = "Yes"

# Row 3: new bug report entry
$ws.Range("A3").Value = "Header Bar Hidding Content"
$ws.Range("B3").Value = "Open Angular application on http://localhost:4200/"
$ws.Range("D3").Value = "Header Blocks first element of Locations List"
$ws.Range("C3").Value = "First Location in list of loactions should be displayed below header nav"
$ws.Range("E3").Value = "Leo"

# Update the active selection to F3
$ws.Range("F3").Select()
